# Update the "想去人数" (interested-count) figures for a handful of events
# that changed, and drop a duplicate row that had crept into the
# "全部类型" (all-types) roll-up sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Bump F (interested-people count) for each event id, everywhere that
#    event shows up (展览 / 演出 / 本地生活 sheets use plain row numbers;
#    全部类型 is handled further down together with the de-dup shift).
# ---------------------------------------------------------------------

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws1.Cells.Item(6, 6).Value2  = 709    # id=80716
$ws1.Cells.Item(7, 6).Value2  = 322    # id=82542
$ws1.Cells.Item(9, 6).Value2  = 118    # id=82486
$ws1.Cells.Item(12, 6).Value2 = 4784   # id=82121
$ws1.Cells.Item(15, 6).Value2 = 467    # id=82319
$ws1.Cells.Item(17, 6).Value2 = 524    # id=81454
$ws1.Cells.Item(18, 6).Value2 = 316    # id=82511
$ws1.Cells.Item(24, 6).Value2 = 289    # id=82056
$ws1.Cells.Item(27, 6).Value2 = 1674   # id=81962
$ws1.Cells.Item(28, 6).Value2 = 400    # id=82458

$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws2.Cells.Item(6, 6).Value2  = 36     # id=82607

$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws3.Cells.Item(3, 6).Value2  = 157    # id=82861

# ---------------------------------------------------------------------
# 2) 全部类型 sheet: apply the same count bumps to the rows that are
#    staying in place (rows 2-16 are untouched by the de-dup below),
#    then remove the duplicated "第九届娃展沙龙" row (old row 17, a
#    copy of row 16) by sliding the B:I content of every following row
#    up by one, and finally deleting the now-spare last row.
# ---------------------------------------------------------------------

$ws4 = $wb.Worksheets.Item(4)   # 全部类型

$ws4.Cells.Item(6, 6).Value2  = 157    # id=82861
$ws4.Cells.Item(10, 6).Value2 = 709    # id=80716
$ws4.Cells.Item(12, 6).Value2 = 322    # id=82542
$ws4.Cells.Item(14, 6).Value2 = 118    # id=82486

# These live in rows that will be slid up by the de-dup shift, so bump
# them on the source row now - the shift below carries the new value.
$ws4.Cells.Item(18, 6).Value2 = 4784   # id=82121
$ws4.Cells.Item(22, 6).Value2 = 467    # id=82319
$ws4.Cells.Item(24, 6).Value2 = 524    # id=81454
$ws4.Cells.Item(25, 6).Value2 = 316    # id=82511
$ws4.Cells.Item(27, 6).Value2 = 36     # id=82607
$ws4.Cells.Item(38, 6).Value2 = 289    # id=82056
$ws4.Cells.Item(41, 6).Value2 = 1674   # id=81962
$ws4.Cells.Item(42, 6).Value2 = 400    # id=82458

# Row 17 is a straight duplicate of row 16 (same event). Drop it by
# copying columns B..I (everything except the sequential index in A)
# from each following row up into the previous row, bottom-to-top is
# not needed here since we always read from the not-yet-overwritten
# source row before writing it.
for ($r = 17; $r -le 42; $r++) {
    for ($c = 2; $c -le 9; $c++) {
        $val = $ws4.Cells.Item($r + 1, $c).Value2
        if ($c -eq 2) {
            # Column B holds plain "YYYY-MM-DD" text. Without the leading
            # apostrophe Excel "helpfully" reinterprets it as a date
            # serial number, which is not what the source file stores.
            $ws4.Cells.Item($r, $c).Value2 = "'" + $val
            $ws4.Cells.Item($r, $c).ClearFormats()
        } else {
            $ws4.Cells.Item($r, $c).Value2 = $val
        }
    }
}

# Row 43 is now a leftover duplicate of the new row 42 - remove it
# outright so the sheet ends at row 42 (dimension A1:I42).
$ws4.Rows("43").Delete()
